# kpi pagination, show user name kpi, alicia maju bisa edit kpi sendiri
#
# The sheet lists daily task rows (9 rows per day) for a week, in column A.
# Shift the week forward by exactly one week: 2023-09-18..23 -> 2023-09-25..30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDates = @("2023-09-18", "2023-09-19", "2023-09-20", "2023-09-21", "2023-09-22", "2023-09-23")
$newDates = @("2023-09-25", "2023-09-26", "2023-09-27", "2023-09-28", "2023-09-29", "2023-09-30")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    for ($i = 0; $i -lt $oldDates.Length; $i++) {
        if ($val -eq $oldDates[$i]) {
            $cell.Value2 = $newDates[$i]
        }
    }
}

# Move the active selection from B47 to D45 (pagination / UI state change).
$ws.Range("D45").Select()
